# TC05_Canine_Filter_Breed-Beagle.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The "CasesTab" query (row 2, column B = "query" column on the "startup"
# sheet) still referenced the retired `co:cohort` pattern / `Cohort` output
# column. That trailing, now-invalid `coalesce(co.cohort_description, '') AS
# `Cohort`` line is removed so the query matches the corrected Cypher used
# everywhere else in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Fix the CasesTab query text (B2): drop the trailing Cohort column ---
$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Beagle']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCasesQuery

# --- Restore the on-disk row heights for the wrapped query cells now that
#     the CasesTab text is one line shorter (matches the re-saved workbook)
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# --- View state: the re-saved workbook scrolls back to the top of the
#     sheet, zooms in to 130%, and leaves the selection on B2 instead of B4
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("B2").Select()
